# Day 13 (column S) attendance marks for participants (rows 7-80) are being
# retracted/cleared — the "P"/"A" marks entered for that day are removed,
# leaving the cell blank (formatted the same way the still-unfilled days,
# i.e. column T onward, already are). The dependent Total-Absence (E) and
# Total-Present (F) COUNTIF formulas recalc automatically once the marks
# are gone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column S (day 13) uses style 29 (matches G:R) while column T (day 14,
# first "blank" day) already carries the right border/format for an empty
# attendance cell. Copy that formatting onto S before clearing it so the
# cleared cells look like every other not-yet-filled-in day, then clear
# the P/A values themselves.
$ws.Range("T7:T80").Copy()
$ws.Range("S7:S80").PasteSpecial(-4122)
$ws.Range("S7:S80").ClearContents()
